$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-150: update date serial from 45186 to 45188
# (Excel serial 45188 corresponds to 2023-09-19)
$ws.Range("C2:C150").Value = 45188
